# Update the "Corr/total marks" figures on the marksheet's scoring summary
# (Marking row and Total row of the quiz worksheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row: number of right-answer marks awarded per question
$ws.Range("B11").Value = 5

# "Total" row: total marks obtained
$ws.Range("B12").Value = 120

# "Total" row, Max column: correct/total marks string
$ws.Range("E12").Value = "120/140"
